# Applies the "Fixed ICDC breed all testcases" edit to the
# TC01_Canine_Filter_Breed-Akita workbook.
#
# Changes:
#  1. Updates the "StatQuery" Cypher text (shared by cells C2:C4 on the
#     "startup" sheet) to the new query text.
#  2. Updates the sheetView of the "startup" sheet: removes the
#     topLeftCell="A4" freeze/scroll position, changes the zoom from 55%
#     to 85%, and moves the active selection from B4 to B2.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

# 1. Update the shared "StatQuery" text used by C2, C3 and C4.
$newQuery = @'
MATCH (demo:demographic)-->(c:case)-->(s:study)-->(p:program)
WHERE demo.breed IN ["Akita"]
OPTIONAL MATCH (c)<-[*]-(samp:sample)
OPTIONAL MATCH (c)<-[*]-(f:file)
OPTIONAL MATCH (f:file)-[*]->(c)
OPTIONAL MATCH (sf:file)-->(s)
  RETURN 
    count(distinct p) AS Programs,
    count(distinct s) AS Studies,
    count(distinct c) AS Cases,
    count(distinct samp) AS Samples,
    count(distinct f) AS `Case Files`,
    count(distinct sf) AS `Study Files`
'@

$ws.Range("C2").Value = $newQuery
$ws.Range("C3").Value = $newQuery
$ws.Range("C4").Value = $newQuery

# 2. Update the sheet view: re-center on A1 (clears topLeftCell), zoom to
#    85%, and move the selection to B2.
$ws.Activate()
$ws.Range("B2").Select()
$excel.ActiveWindow.Zoom = 85
$excel.ActiveWindow.ZoomScaleNormal = 85
